# Apply the regenerated linear-experiment values (ex9.1.10) to the workbook.
# The underlying data in the workbook is stored as TEXT (even values that look
# like numbers), so we force a Text number format before writing the values
# and clear the format afterwards so the cells keep the default style.

$wb = $excel.ActiveWorkbook

# NOTE: worksheet lookup by name is case-insensitive in this environment, and
# this workbook has two sheets whose names only differ by case
# ("Vector_bf" vs "Vector_BF"), so we must address sheets by their (1-based)
# index to avoid ambiguity.
$wsFollower = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$wsPunto    = $wb.Worksheets.Item(4)   # Punto_modificado
$wsVecbf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBF    = $wb.Worksheets.Item(6)   # Vector_BF

# Make sure the ranges we are about to touch are forced to Text so that the
# numeric-looking strings are not silently re-interpreted as numbers.
$wsFollower.Range("A2:F5").NumberFormat = "@"
$wsPunto.Range("A2:D2").NumberFormat = "@"
$wsVecbf.Range("A2:A3").NumberFormat = "@"
$wsVecBF.Range("A2:A5").NumberFormat = "@"

# --- Restricciones_del_follower ---------------------------------------
$wsFollower.Range("A2").Value = "1.7500000000000169 - 2x_1 + y_1 - y_2"
$wsFollower.Range("B2").Value = "0.7499999999999831"
$wsFollower.Range("D2").Value = "0.77"
$wsFollower.Range("E2").Value = "0"
$wsFollower.Range("F2").Value = "7.9"

$wsFollower.Range("A3").Value = "6.549999999999994 + x_1 - 3x_2 + y_2"
$wsFollower.Range("B3").Value = "-8.549999999999994"
$wsFollower.Range("D3").Value = "0.46"
$wsFollower.Range("E3").Value = "4.5"
$wsFollower.Range("F3").Value = "0"

$wsFollower.Range("A4").Value = "104.6 - y_1"
$wsFollower.Range("B4").Value = "-104.6"
$wsFollower.Range("D4").Value = "0.41"
$wsFollower.Range("E4").Value = "0"
$wsFollower.Range("F4").Value = "6.800000000000001"

$wsFollower.Range("A5").Value = "-2.05 - y_2"
$wsFollower.Range("B5").Value = "-2.05"
$wsFollower.Range("D5").Value = "0.64"
$wsFollower.Range("E5").Value = "0"
$wsFollower.Range("F5").Value = "7.5"

# --- Punto_modificado ---------------------------------------------------
$wsPunto.Range("A2").Value = "52.150000000000006"
$wsPunto.Range("B2").Value = "20.25"
$wsPunto.Range("C2").Value = "104.6"
$wsPunto.Range("D2").Value = "2.05"

# --- Vector_bf ------------------------------------------------------------
$wsVecbf.Range("A2").Value = "3.64"
$wsVecbf.Range("A3").Value = "-0.050000000000000044"

# --- Vector_BF ------------------------------------------------------------
$wsVecBF.Range("A2").Value = "-2.5"
$wsVecBF.Range("A3").Value = "12.5"
$wsVecBF.Range("A4").Value = "-0.5"
$wsVecBF.Range("A5").Value = "-4.5"

# Remove the temporary Text formatting so the cells fall back to the default
# (unformatted) style, matching the original workbook's styling.
$wsFollower.Range("A2:F5").ClearFormats()
$wsPunto.Range("A2:D2").ClearFormats()
$wsVecbf.Range("A2:A3").ClearFormats()
$wsVecBF.Range("A2:A5").ClearFormats()
